$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.859.31'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").Value = '2.607.83'
$ws.Range("E3").Value = '  +1.17%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'578.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.16%  '
$ws.Range("D6").Value = "'143.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("E7").Value = '  -0.26%  '
$ws.Range("E8").Value = '  +0.50%  '
$ws.Range("D9").Value = '2.634.42'
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("D10").Value = "'6.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.98%  '
$ws.Range("E11").Value = '  +2.31%  '
$ws.Range("E12").Value = '  -5.21%  '
$ws.Range("E13").Value = '  +5.64%  '
$ws.Range("D14").Value = '3.075.22'
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("D15").Value = '60.842.91'
$ws.Range("E15").Value = '  +2.56%  '
$ws.Range("D16").Value = "'23.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.81%  '
$ws.Range("E17").Value = '  +4.47%  '
$ws.Range("D18").Value = '2.626.48'
$ws.Range("E18").Value = '  +1.75%  '
$ws.Range("D19").Value = "'11.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.25%  '
$ws.Range("E20").Value = '  +2.79%  '
$ws.Range("D21").Value = "'350.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.79%  '
$ws.Range("E22").Value = '  +7.33%  '
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = "'0.517"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.47%  '
$ws.Range("D25").Value = "'63.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("D28").Value = "'7.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.10%  '
$ws.Range("E30").Value = '  +8.86%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = "'6.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.44%  '
$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D32").Value = "'0.997"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("D33").Value = "'162.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.17%  '
$ws.Range("D34").Value = "'19.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.65%  '
$ws.Range("E35").Value = '  +13.57%  '
$ws.Range("E36").Value = '  +4.70%  '
$ws.Range("D38").Value = "'1.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.52%  '
$ws.Range("D39").Value = "'37.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.46%  '
$ws.Range("E40").Value = '  +5.74%  '
$ws.Range("D41").Value = "'309.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.72%  '
$ws.Range("D42").Value = "'0.848"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.42%  '
$ws.Range("D43").Value = "'133.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.36%  '
$ws.Range("D44").Value = "'20.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.97%  '
$ws.Range("D45").Value = "'19.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.77%  '
$ws.Range("D46").Value = "'5.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +12.09%  '
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = "'0.610"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.80%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = "'0.0985"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.26%  '
$ws.Range("D50").Value = "'0.0552"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.00%  '
$ws.Range("E51").Value = '  +3.60%  '
